$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 12 blank "B" condition rows (old rows 26-37)
$ws.Rows("26:37").Delete()

# Rewrite the stimuli table with the new run (trialnum 193-240)
$ws.Cells.Item(2,1).Value = "L.png"
$ws.Cells.Item(2,2).Value = 193
$ws.Cells.Item(2,3).Value = "R"
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(3,1).Value = "M.png"
$ws.Cells.Item(3,2).Value = 194
$ws.Cells.Item(3,3).Value = "R"
$ws.Cells.Item(3,4).Value = 0
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(4,1).Value = "E.png"
$ws.Cells.Item(4,2).Value = 195
$ws.Cells.Item(4,3).Value = "R"
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(5,1).Value = "F.png"
$ws.Cells.Item(5,2).Value = 196
$ws.Cells.Item(5,3).Value = "R"
$ws.Cells.Item(5,4).Value = 0
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(6,1).Value = "E.png"
$ws.Cells.Item(6,2).Value = 197
$ws.Cells.Item(6,3).Value = "R"
$ws.Cells.Item(6,4).Value = 0
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(7,1).Value = "D.png"
$ws.Cells.Item(7,2).Value = 198
$ws.Cells.Item(7,3).Value = "R"
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(8,1).Value = "C.png"
$ws.Cells.Item(8,2).Value = 199
$ws.Cells.Item(8,3).Value = "R"
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(9,1).Value = "L.png"
$ws.Cells.Item(9,2).Value = 200
$ws.Cells.Item(9,3).Value = "R"
$ws.Cells.Item(9,4).Value = 0
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(10,1).Value = "A.png"
$ws.Cells.Item(10,2).Value = 201
$ws.Cells.Item(10,3).Value = "R"
$ws.Cells.Item(10,4).Value = 0
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(11,1).Value = "B.png"
$ws.Cells.Item(11,2).Value = 202
$ws.Cells.Item(11,3).Value = "R"
$ws.Cells.Item(11,4).Value = 0
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(12,1).Value = "B.png"
$ws.Cells.Item(12,2).Value = 203
$ws.Cells.Item(12,3).Value = "R"
$ws.Cells.Item(12,4).Value = 0
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(13,1).Value = "H.png"
$ws.Cells.Item(13,2).Value = 204
$ws.Cells.Item(13,3).Value = "R"
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(14,1).Value = "G.png"
$ws.Cells.Item(14,2).Value = 205
$ws.Cells.Item(14,3).Value = "R"
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).Value = 1
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(15,1).Value = "K.png"
$ws.Cells.Item(15,2).Value = 206
$ws.Cells.Item(15,3).Value = "R"
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 2
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(16,1).Value = "G.png"
$ws.Cells.Item(16,2).Value = 207
$ws.Cells.Item(16,3).Value = "R"
$ws.Cells.Item(16,4).Value = 0
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(17,1).Value = "A.png"
$ws.Cells.Item(17,2).Value = 208
$ws.Cells.Item(17,3).Value = "R"
$ws.Cells.Item(17,4).Value = 0
$ws.Cells.Item(17,5).Value = 1
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(18,1).Value = "M.png"
$ws.Cells.Item(18,2).Value = 209
$ws.Cells.Item(18,3).Value = "R"
$ws.Cells.Item(18,4).Value = 0
$ws.Cells.Item(18,5).Value = 2
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(19,1).Value = "K.png"
$ws.Cells.Item(19,2).Value = 210
$ws.Cells.Item(19,3).Value = "R"
$ws.Cells.Item(19,4).Value = 0
$ws.Cells.Item(19,5).Value = 3
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(20,1).Value = "C.png"
$ws.Cells.Item(20,2).Value = 211
$ws.Cells.Item(20,3).Value = "R"
$ws.Cells.Item(20,4).Value = 0
$ws.Cells.Item(20,5).Value = 1
$ws.Cells.Item(20,6).Value = 1
$ws.Cells.Item(21,1).Value = "J.png"
$ws.Cells.Item(21,2).Value = 212
$ws.Cells.Item(21,3).Value = "R"
$ws.Cells.Item(21,4).Value = 0
$ws.Cells.Item(21,5).Value = 2
$ws.Cells.Item(21,6).Value = 1
$ws.Cells.Item(22,1).Value = "J.png"
$ws.Cells.Item(22,2).Value = 213
$ws.Cells.Item(22,3).Value = "R"
$ws.Cells.Item(22,4).Value = 0
$ws.Cells.Item(22,5).Value = 3
$ws.Cells.Item(22,6).Value = 1
$ws.Cells.Item(23,1).Value = "H.png"
$ws.Cells.Item(23,2).Value = 214
$ws.Cells.Item(23,3).Value = "R"
$ws.Cells.Item(23,4).Value = 0
$ws.Cells.Item(23,5).Value = 1
$ws.Cells.Item(23,6).Value = 1
$ws.Cells.Item(24,1).Value = "D.png"
$ws.Cells.Item(24,2).Value = 215
$ws.Cells.Item(24,3).Value = "R"
$ws.Cells.Item(24,4).Value = 0
$ws.Cells.Item(24,5).Value = 2
$ws.Cells.Item(24,6).Value = 1
$ws.Cells.Item(25,1).Value = "F.png"
$ws.Cells.Item(25,2).Value = 216
$ws.Cells.Item(25,3).Value = "R"
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(25,5).Value = 3
$ws.Cells.Item(25,6).Value = 1
$ws.Cells.Item(26,1).Value = "F.png"
$ws.Cells.Item(26,2).Value = 217
$ws.Cells.Item(26,3).Value = "R"
$ws.Cells.Item(26,4).Value = 0
$ws.Cells.Item(26,5).Value = 1
$ws.Cells.Item(26,6).Value = 1
$ws.Cells.Item(27,1).Value = "K.png"
$ws.Cells.Item(27,2).Value = 218
$ws.Cells.Item(27,3).Value = "R"
$ws.Cells.Item(27,4).Value = 0
$ws.Cells.Item(27,5).Value = 2
$ws.Cells.Item(27,6).Value = 1
$ws.Cells.Item(28,1).Value = "C.png"
$ws.Cells.Item(28,2).Value = 219
$ws.Cells.Item(28,3).Value = "R"
$ws.Cells.Item(28,4).Value = 0
$ws.Cells.Item(28,5).Value = 3
$ws.Cells.Item(28,6).Value = 1
$ws.Cells.Item(29,1).Value = "G.png"
$ws.Cells.Item(29,2).Value = 220
$ws.Cells.Item(29,3).Value = "R"
$ws.Cells.Item(29,4).Value = 0
$ws.Cells.Item(29,5).Value = 1
$ws.Cells.Item(29,6).Value = 1
$ws.Cells.Item(30,1).Value = "A.png"
$ws.Cells.Item(30,2).Value = 221
$ws.Cells.Item(30,3).Value = "R"
$ws.Cells.Item(30,4).Value = 0
$ws.Cells.Item(30,5).Value = 2
$ws.Cells.Item(30,6).Value = 1
$ws.Cells.Item(31,1).Value = "M.png"
$ws.Cells.Item(31,2).Value = 222
$ws.Cells.Item(31,3).Value = "R"
$ws.Cells.Item(31,4).Value = 0
$ws.Cells.Item(31,5).Value = 3
$ws.Cells.Item(31,6).Value = 1
$ws.Cells.Item(32,1).Value = "C.png"
$ws.Cells.Item(32,2).Value = 223
$ws.Cells.Item(32,3).Value = "R"
$ws.Cells.Item(32,4).Value = 0
$ws.Cells.Item(32,5).Value = 1
$ws.Cells.Item(32,6).Value = 1
$ws.Cells.Item(33,1).Value = "A.png"
$ws.Cells.Item(33,2).Value = 224
$ws.Cells.Item(33,3).Value = "R"
$ws.Cells.Item(33,4).Value = 0
$ws.Cells.Item(33,5).Value = 2
$ws.Cells.Item(33,6).Value = 1
$ws.Cells.Item(34,1).Value = "J.png"
$ws.Cells.Item(34,2).Value = 225
$ws.Cells.Item(34,3).Value = "R"
$ws.Cells.Item(34,4).Value = 0
$ws.Cells.Item(34,5).Value = 3
$ws.Cells.Item(34,6).Value = 1
$ws.Cells.Item(35,1).Value = "E.png"
$ws.Cells.Item(35,2).Value = 226
$ws.Cells.Item(35,3).Value = "R"
$ws.Cells.Item(35,4).Value = 0
$ws.Cells.Item(35,5).Value = 1
$ws.Cells.Item(35,6).Value = 1
$ws.Cells.Item(36,1).Value = "J.png"
$ws.Cells.Item(36,2).Value = 227
$ws.Cells.Item(36,3).Value = "R"
$ws.Cells.Item(36,4).Value = 0
$ws.Cells.Item(36,5).Value = 2
$ws.Cells.Item(36,6).Value = 1
$ws.Cells.Item(37,1).Value = "H.png"
$ws.Cells.Item(37,2).Value = 228
$ws.Cells.Item(37,3).Value = "R"
$ws.Cells.Item(37,4).Value = 0
$ws.Cells.Item(37,5).Value = 3
$ws.Cells.Item(37,6).Value = 1
$ws.Cells.Item(38,1).Value = "H.png"
$ws.Cells.Item(38,2).Value = 229
$ws.Cells.Item(38,3).Value = "R"
$ws.Cells.Item(38,4).Value = 0
$ws.Cells.Item(38,5).Value = 1
$ws.Cells.Item(38,6).Value = 1
$ws.Cells.Item(39,1).Value = "M.png"
$ws.Cells.Item(39,2).Value = 230
$ws.Cells.Item(39,3).Value = "R"
$ws.Cells.Item(39,4).Value = 0
$ws.Cells.Item(39,5).Value = 2
$ws.Cells.Item(39,6).Value = 1
$ws.Cells.Item(40,1).Value = "D.png"
$ws.Cells.Item(40,2).Value = 231
$ws.Cells.Item(40,3).Value = "R"
$ws.Cells.Item(40,4).Value = 0
$ws.Cells.Item(40,5).Value = 3
$ws.Cells.Item(40,6).Value = 1
$ws.Cells.Item(41,1).Value = "D.png"
$ws.Cells.Item(41,2).Value = 232
$ws.Cells.Item(41,3).Value = "R"
$ws.Cells.Item(41,4).Value = 0
$ws.Cells.Item(41,5).Value = 1
$ws.Cells.Item(41,6).Value = 1
$ws.Cells.Item(42,1).Value = "L.png"
$ws.Cells.Item(42,2).Value = 233
$ws.Cells.Item(42,3).Value = "R"
$ws.Cells.Item(42,4).Value = 0
$ws.Cells.Item(42,5).Value = 2
$ws.Cells.Item(42,6).Value = 1
$ws.Cells.Item(43,1).Value = "B.png"
$ws.Cells.Item(43,2).Value = 234
$ws.Cells.Item(43,3).Value = "R"
$ws.Cells.Item(43,4).Value = 0
$ws.Cells.Item(43,5).Value = 3
$ws.Cells.Item(43,6).Value = 1
$ws.Cells.Item(44,1).Value = "K.png"
$ws.Cells.Item(44,2).Value = 235
$ws.Cells.Item(44,3).Value = "R"
$ws.Cells.Item(44,4).Value = 0
$ws.Cells.Item(44,5).Value = 1
$ws.Cells.Item(44,6).Value = 1
$ws.Cells.Item(45,1).Value = "E.png"
$ws.Cells.Item(45,2).Value = 236
$ws.Cells.Item(45,3).Value = "R"
$ws.Cells.Item(45,4).Value = 0
$ws.Cells.Item(45,5).Value = 2
$ws.Cells.Item(45,6).Value = 1
$ws.Cells.Item(46,1).Value = "F.png"
$ws.Cells.Item(46,2).Value = 237
$ws.Cells.Item(46,3).Value = "R"
$ws.Cells.Item(46,4).Value = 0
$ws.Cells.Item(46,5).Value = 3
$ws.Cells.Item(46,6).Value = 1
$ws.Cells.Item(47,1).Value = "G.png"
$ws.Cells.Item(47,2).Value = 238
$ws.Cells.Item(47,3).Value = "R"
$ws.Cells.Item(47,4).Value = 0
$ws.Cells.Item(47,5).Value = 1
$ws.Cells.Item(47,6).Value = 1
$ws.Cells.Item(48,1).Value = "B.png"
$ws.Cells.Item(48,2).Value = 239
$ws.Cells.Item(48,3).Value = "R"
$ws.Cells.Item(48,4).Value = 0
$ws.Cells.Item(48,5).Value = 2
$ws.Cells.Item(48,6).Value = 1
$ws.Cells.Item(49,1).Value = "L.png"
$ws.Cells.Item(49,2).Value = 240
$ws.Cells.Item(49,3).Value = "R"
$ws.Cells.Item(49,4).Value = 0
$ws.Cells.Item(49,5).Value = 3
$ws.Cells.Item(49,6).Value = 1

# Update the view: drop the old scroll position, select K25 like the source file
$ws.Range("K25").Select()
